# feat: changed fecha pago
# The "WarehouseCode"/"WhsCode" column (column F) on Sheet1 is a duplicate
# of the existing WhsCode column and is removed entirely, shifting all
# subsequent columns (G:N) left by one (to F:M).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Columns.Item(6).Delete()
